$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values for the interpolation example
$ws.Range("A7").Value = -29.628599999999999
$ws.Range("B7").Value = 5163.7
$ws.Range("D7").Value = 0
$ws.Range("A9").Value = 30.505700000000001
$ws.Range("B9").Value = 5164.3

# Move the active selection to F9, matching the author's final cursor position
$ws.Range("F9").Select()
